$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report date range) ---
$ws.Range("A8").Value = "Volume 31   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  9/30/2024  Through  10/6/2024"

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -75.471698113207
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 34
$ws.Range("K15").Value = -2.941176470588
$ws.Range("L15").Value = 6.451612903225
$ws.Range("N15").Value = -53.521126760563
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = 18.181818181818
$ws.Range("F16").Value = 48
$ws.Range("G16").Value = 31
$ws.Range("H16").Value = 54.838709677419
$ws.Range("I16").Value = 434
$ws.Range("J16").Value = 426
$ws.Range("K16").Value = 1.8779342723
$ws.Range("L16").Value = -20.220588235294
$ws.Range("M16").Value = 17.934782608695
$ws.Range("N16").Value = -71.105193075898
$ws.Range("C17").Value = 24
$ws.Range("D17").Value = 23
$ws.Range("E17").Value = 4.347826086956
$ws.Range("F17").Value = 87
$ws.Range("G17").Value = 83
$ws.Range("H17").Value = 4.819277108433
$ws.Range("I17").Value = 815
$ws.Range("J17").Value = 862
$ws.Range("K17").Value = -5.452436194895
$ws.Range("L17").Value = 0.246002460024
$ws.Range("M17").Value = 89.53488372093
$ws.Range("N17").Value = -24.537037037037
$ws.Range("C18").Value = 9
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = -10
$ws.Range("F18").Value = 31
$ws.Range("G18").Value = 29
$ws.Range("H18").Value = 6.896551724137
$ws.Range("I18").Value = 224
$ws.Range("J18").Value = 368
$ws.Range("K18").Value = -39.130434782608
$ws.Range("L18").Value = -31.288343558282
$ws.Range("M18").Value = 10.89108910891
$ws.Range("N18").Value = -86.112833230006
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -12.5
$ws.Range("F19").Value = 67
$ws.Range("G19").Value = 72
$ws.Range("H19").Value = -6.944444444444
$ws.Range("I19").Value = 663
$ws.Range("J19").Value = 619
$ws.Range("K19").Value = 7.108239095315
$ws.Range("L19").Value = -10.887096774193
$ws.Range("M19").Value = 123.986486486486
$ws.Range("N19").Value = 5.071315372424
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -33.333333333333
$ws.Range("G20").Value = 33
$ws.Range("H20").Value = -12.121212121212
$ws.Range("I20").Value = 281
$ws.Range("J20").Value = 373
$ws.Range("K20").Value = -24.664879356568
$ws.Range("L20").Value = -5.704697986577
$ws.Range("M20").Value = 70.30303030303
$ws.Range("N20").Value = -72.797676669893
$ws.Range("C21").Value = 64
$ws.Range("D21").Value = 67
$ws.Range("E21").Value = -4.477611940298
$ws.Range("F21").Value = 264
$ws.Range("G21").Value = 251
$ws.Range("H21").Value = 5.179282868525
$ws.Range("I21").Value = 2463
$ws.Range("J21").Value = 2695
$ws.Range("K21").Value = -8.60853432282
$ws.Range("L21").Value = -11.211247296323
$ws.Range("M21").Value = 64.419225634178
$ws.Range("N21").Value = -58.833361190038
$ws.Range("L22").Value = 36.666666666666
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = 100
$ws.Range("I23").Value = 62
$ws.Range("J23").Value = 63
$ws.Range("K23").Value = -1.587301587301
$ws.Range("L23").Value = 3.333333333333
$ws.Range("M23").Value = 63.157894736842
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 47
$ws.Range("E24").Value = -42.553191489361
$ws.Range("F24").Value = 164
$ws.Range("G24").Value = 186
$ws.Range("H24").Value = -11.827956989247
$ws.Range("I24").Value = 1510
$ws.Range("J24").Value = 1620
$ws.Range("K24").Value = -6.79012345679
$ws.Range("L24").Value = 1.889338731443
$ws.Range("M24").Value = 41.385767790262
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 28
$ws.Range("E25").Value = -75
$ws.Range("F25").Value = 65
$ws.Range("G25").Value = 104
$ws.Range("H25").Value = -37.5
$ws.Range("I25").Value = 700
$ws.Range("J25").Value = 845
$ws.Range("K25").Value = -17.159763313609
$ws.Range("L25").Value = -4.371584699453
$ws.Range("C26").Value = 35
$ws.Range("D26").Value = 26
$ws.Range("E26").Value = 34.615384615384
$ws.Range("F26").Value = 138
$ws.Range("G26").Value = 81
$ws.Range("H26").Value = 70.37037037037
$ws.Range("I26").Value = 1121
$ws.Range("J26").Value = 953
$ws.Range("K26").Value = 17.628541448058
$ws.Range("L26").Value = 11.320754716981
$ws.Range("M26").Value = 15.092402464065
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = -66.666666666666
$ws.Range("J27").Value = 57
$ws.Range("K27").Value = -8.771929824561
$ws.Range("L27").Value = 6.122448979591
$ws.Range("C28").Value = 3
$ws.Range("F28").Value = 11
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 102
$ws.Range("J28").Value = 106
$ws.Range("K28").Value = -3.77358490566
$ws.Range("L28").Value = 12.087912087912
$ws.Range("F29").Value = 6
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 49
$ws.Range("K29").Value = 10.204081632653
$ws.Range("M29").Value = 54.285714285714
$ws.Range("N29").Value = -62.5
$ws.Range("F30").Value = 4
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 39
$ws.Range("K30").Value = 0
$ws.Range("M30").Value = 39.285714285714
$ws.Range("N30").Value = -69.53125
$ws.Range("F33").Value = 1

# --- Numeric -> text marker conversions (style must follow a text-styled donor cell) ---
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("D22").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "0"
$ws.Range("D22").Copy()
$ws.Range("G22").PasteSpecial(-4122)

$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "***.*"
$ws.Range("E22").Copy()
$ws.Range("H22").PasteSpecial(-4122)

# --- Text marker -> numeric conversions (style must follow a numeric-styled donor cell) ---
$ws.Range("D28").Value = 3
$ws.Range("C28").Copy()
$ws.Range("D28").PasteSpecial(-4122)

$ws.Range("E28").Value = 0
$ws.Range("H28").Copy()
$ws.Range("E28").PasteSpecial(-4122)

$ws.Range("D29").Value = 1
$ws.Range("C28").Copy()
$ws.Range("D29").PasteSpecial(-4122)

$ws.Range("E29").Value = -100
$ws.Range("H28").Copy()
$ws.Range("E29").PasteSpecial(-4122)

$ws.Range("D30").Value = 1
$ws.Range("C28").Copy()
$ws.Range("D30").PasteSpecial(-4122)

$ws.Range("E30").Value = -100
$ws.Range("H28").Copy()
$ws.Range("E30").PasteSpecial(-4122)

$excel.CutCopyMode = $false

